$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated strikeout ("K") values for rows 2-12 (column G)
$kValues = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 3
    6  = 1
    7  = 4
    8  = 1
    9  = 1
    10 = 0
    11 = 2
    12 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
